$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
foreach ($addr in @("D2","E2")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D2').Value = '27.353.80'
$ws.Range('E2').Value = '  +6.90%  '
foreach ($addr in @("D2","E2")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 3
foreach ($addr in @("D3","E3")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D3').Value = '1.810.63'
$ws.Range('E3').Value = '  +6.25%  '
foreach ($addr in @("D3","E3")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 4
foreach ($addr in @("D4","E4")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.73%  '
foreach ($addr in @("D4","E4")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 5
foreach ($addr in @("D5","E5")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D5').Value = '343.05'
$ws.Range('E5').Value = '  +4.18%  '
foreach ($addr in @("D5","E5")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 6
foreach ($addr in @("D6","E6")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.28%  '
foreach ($addr in @("D6","E6")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 7
foreach ($addr in @("D7","E7")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D7').Value = '0.3828'
$ws.Range('E7').Value = '  +4.34%  '
foreach ($addr in @("D7","E7")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 8
foreach ($addr in @("D8","E8")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D8').Value = '50.31'
$ws.Range('E8').Value = '  +5.42%  '
foreach ($addr in @("D8","E8")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 9
foreach ($addr in @("D9","E9")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D9').Value = '0.3515'
$ws.Range('E9').Value = '  +7.17%  '
foreach ($addr in @("D9","E9")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 10
foreach ($addr in @("D10","E10")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D10').Value = '1.229'
$ws.Range('E10').Value = '  +5.72%  '
foreach ($addr in @("D10","E10")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 11
foreach ($addr in @("D11","E11")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D11').Value = '0.07789'
$ws.Range('E11').Value = '  +6.42%  '
foreach ($addr in @("D11","E11")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 12
foreach ($addr in @("D12","E12")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D12').Value = '1.005'
$ws.Range('E12').Value = '  +0.89%  '
foreach ($addr in @("D12","E12")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 13
foreach ($addr in @("D13","E13")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D13').Value = '22.33'
$ws.Range('E13').Value = '  +12.22%  '
foreach ($addr in @("D13","E13")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 14
foreach ($addr in @("D14","E14")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D14').Value = '6.628'
$ws.Range('E14').Value = '  +7.25%  '
foreach ($addr in @("D14","E14")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 15
foreach ($addr in @("D15","E15")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D15').Value = '7.234'
$ws.Range('E15').Value = '  +6.65%  '
foreach ($addr in @("D15","E15")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 16
foreach ($addr in @("D16","E16")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D16').Value = '1.814.70'
$ws.Range('E16').Value = '  +6.82%  '
foreach ($addr in @("D16","E16")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 17
foreach ($addr in @("D17","E17")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D17').Value = '0.00001123'
$ws.Range('E17').Value = '  +5.21%  '
foreach ($addr in @("D17","E17")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 18
foreach ($addr in @("D18","E18")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D18').Value = '0.06761'
$ws.Range('E18').Value = '  +2.83%  '
foreach ($addr in @("D18","E18")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 19
foreach ($addr in @("D19","E19")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D19').Value = '87.18'
$ws.Range('E19').Value = '  +8.22%  '
foreach ($addr in @("D19","E19")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 20
foreach ($addr in @("D20","E20")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.43%  '
foreach ($addr in @("D20","E20")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 21
foreach ($addr in @("D21","E21")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D21').Value = '17.87'
$ws.Range('E21').Value = '  +11.27%  '
foreach ($addr in @("D21","E21")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 22
foreach ($addr in @("D22","E22")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D22').Value = '6.545'
$ws.Range('E22').Value = '  +8.57%  '
foreach ($addr in @("D22","E22")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 23
foreach ($addr in @("D23","E23")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D23').Value = '13.11'
$ws.Range('E23').Value = '  +0.42%  '
foreach ($addr in @("D23","E23")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 24
foreach ($addr in @("D24","E24")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D24').Value = '27.409.06'
$ws.Range('E24').Value = '  +7.27%  '
foreach ($addr in @("D24","E24")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 25
foreach ($addr in @("D25","E25")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D25').Value = '2.473'
$ws.Range('E25').Value = '  +1.32%  '
foreach ($addr in @("D25","E25")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 26
foreach ($addr in @("D26","E26")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D26').Value = '2.683'
$ws.Range('E26').Value = '  +8.72%  '
foreach ($addr in @("D26","E26")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 27
foreach ($addr in @("D27","E27")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D27').Value = '22.01'
$ws.Range('E27').Value = '  +15.40%  '
foreach ($addr in @("D27","E27")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 28
foreach ($addr in @("D28","E28")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D28').Value = '1.510'
$ws.Range('E28').Value = '  +18.59%  '
foreach ($addr in @("D28","E28")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 29
foreach ($addr in @("D29","E29")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D29').Value = '153.16'
$ws.Range('E29').Value = '  +2.33%  '
foreach ($addr in @("D29","E29")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 30
foreach ($addr in @("D30","E30")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D30').Value = '2.019.93'
$ws.Range('E30').Value = '  +6.94%  '
foreach ($addr in @("D30","E30")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 31
foreach ($addr in @("D31","E31")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D31').Value = '137.30'
$ws.Range('E31').Value = '  +7.24%  '
foreach ($addr in @("D31","E31")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 32
foreach ($addr in @("D32","E32")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D32').Value = '6.396'
$ws.Range('E32').Value = '  +7.60%  '
foreach ($addr in @("D32","E32")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 33
foreach ($addr in @("D33","E33")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D33').Value = '4.163'
$ws.Range('E33').Value = '  +1.60%  '
foreach ($addr in @("D33","E33")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 34
foreach ($addr in @("D34","E34")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D34').Value = '13.90'
$ws.Range('E34').Value = '  +10.09%  '
foreach ($addr in @("D34","E34")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 35
foreach ($addr in @("D35","E35")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D35').Value = '0.08784'
$ws.Range('E35').Value = '  +3.92%  '
foreach ($addr in @("D35","E35")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 36
foreach ($addr in @("E36")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('E36').Value = '  +2.84%  '
foreach ($addr in @("E36")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 37
foreach ($addr in @("D37","E37")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D37').Value = '5.662'
$ws.Range('E37').Value = '  +7.42%  '
foreach ($addr in @("D37","E37")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 38
foreach ($addr in @("D38","E38")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D38').Value = '0.06540'
$ws.Range('E38').Value = '  +5.59%  '
foreach ($addr in @("D38","E38")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 39
foreach ($addr in @("B39","C39","D39","E39")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.6920'
$ws.Range('E39').Value = '  +14.18%  '
foreach ($addr in @("B39","C39","D39","E39")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 40
foreach ($addr in @("D40","E40")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D40').Value = '0.2264'
$ws.Range('E40').Value = '  +7.17%  '
foreach ($addr in @("D40","E40")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 41
foreach ($addr in @("B41","C41","D41","E41")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.02415'
$ws.Range('E41').Value = '  +6.72%  '
foreach ($addr in @("B41","C41","D41","E41")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 42
foreach ($addr in @("D42","E42")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D42').Value = '8.988'
$ws.Range('E42').Value = '  +6.33%  '
foreach ($addr in @("D42","E42")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 43
foreach ($addr in @("D43","E43")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D43').Value = '1.250'
$ws.Range('E43').Value = '  -0.99%  '
foreach ($addr in @("D43","E43")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 44
foreach ($addr in @("D44","E44")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D44').Value = '14.97'
$ws.Range('E44').Value = '  +7.88%  '
foreach ($addr in @("D44","E44")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 45
foreach ($addr in @("D45","E45")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D45').Value = '0.6522'
$ws.Range('E45').Value = '  +11.99%  '
foreach ($addr in @("D45","E45")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 46
foreach ($addr in @("D46","E46")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.43%  '
foreach ($addr in @("D46","E46")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 47
foreach ($addr in @("D47","E47")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D47').Value = '4.023'
$ws.Range('E47').Value = '  +5.04%  '
foreach ($addr in @("D47","E47")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 48
foreach ($addr in @("D48","E48")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D48').Value = '2.171'
$ws.Range('E48').Value = '  +9.01%  '
foreach ($addr in @("D48","E48")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 49
foreach ($addr in @("D49","E49")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D49').Value = '132.79'
$ws.Range('E49').Value = '  +5.94%  '
foreach ($addr in @("D49","E49")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 50
foreach ($addr in @("D50","E50")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D50').Value = '0.07354'
$ws.Range('E50').Value = '  +2.37%  '
foreach ($addr in @("D50","E50")) {
    $ws.Range($addr).Style = "Normal"
}

# Row 51
foreach ($addr in @("D51","E51")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D51').Value = '80.76'
$ws.Range('E51').Value = '  +6.65%  '
foreach ($addr in @("D51","E51")) {
    $ws.Range($addr).Style = "Normal"
}
